# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.474.08"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.551.41"
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.65"
$ws.Range("E5").Value = "  -1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.482"
$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.13"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("E9").Value = "  -1.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0583"
$ws.Range("E10").Value = "  -1.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0890"
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.770.61"
$ws.Range("E12").Value = "  -1.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.549.37"
$ws.Range("E13").Value = "  -2.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.414.92"
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("E15").Value = "  -2.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.509"
$ws.Range("E16").Value = "  -2.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.89"
$ws.Range("E17").Value = "  -2.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.90"
$ws.Range("E18").Value = "  -0.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.33"
$ws.Range("E19").Value = "  -1.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0671"
$ws.Range("E20").Value = "  -2.86%  "

$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.89"
$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.92"
$ws.Range("E23").Value = "  -2.62%  "

$ws.Range("E24").Value = "  -2.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.13"
$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.79"
$ws.Range("E26").Value = "  -1.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.103"
$ws.Range("E27").Value = "  -1.43%  "

$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("E29").Value = "  -3.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0467"
$ws.Range("E30").Value = "  -3.40%  "

$ws.Range("E31").Value = "  -4.85%  "

$ws.Range("E32").Value = "  -1.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.383.50"
$ws.Range("E33").Value = "  -1.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("E34").Value = "  -3.61%  "

$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("E36").Value = "  -3.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").Value = "  -2.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.57"
$ws.Range("E38").Value = "  -3.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0162"
$ws.Range("E39").Value = "  -2.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("E40").Value = "  +1.78%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("B42").Value = "ImmutableX"
$ws.Range("C42").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.508"
$ws.Range("E42").Value = "  -2.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.770"
$ws.Range("E43").Value = "  -2.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0458"
$ws.Range("E44").Value = "  -1.18%  "

$ws.Range("E45").Value = "  -2.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.81"
$ws.Range("E46").Value = "  -2.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.682.98"
$ws.Range("E47").Value = "  -2.16%  "

$ws.Range("E48").Value = "  -9.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.34"
$ws.Range("E49").Value = "  +11.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.73"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  +5.00%  "
